$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad" / Changed date) holds the same serial date value (45180)
# for every data row (2 through 141). Bump it by one day to 45181, matching the
# diff which updates every C2:C141 cell from 45180 to 45181.
for ($row = 2; $row -le 141; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45180) {
        $cell.Value2 = 45181
    }
}
